# below_ground.xlsx edit script
# Implements:
#  - bump sheetId of "fixed" sheet from 4 to 5 (by recreating it in place)
#  - populate "fixed" sheet data (rows 1-4 preserved, rows 5-6 new)
#  - adjust "fixed" sheet column widths
#  - change level1!I27 from 0.1 to 0
#  - update sheet view selections / active tab across sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Recreate the "fixed" sheet so it gets a new (higher) sheetId,
#    while keeping the same tab position and the same r:id.
#    Copying (rather than adding a blank sheet) preserves the
#    original data, column widths and namespace declarations.
# ---------------------------------------------------------------
$oldFixed = $wb.Worksheets.Item("fixed")
$oldFixed.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$oldFixed2 = $wb.Worksheets.Item("fixed")
$oldFixed2.Delete()
$ws4 = $wb.Worksheets.Item("fixed (2)")
$ws4.Name = "fixed"

# ---------------------------------------------------------------
# 3. Add the two new rows (order chosen to reproduce shared-string
#    table layout of the target workbook).
# ---------------------------------------------------------------
$ws4.Range("A5").Value = "source_for_landslide_deformation_geometry"
$ws4.Range("A6").Value = "source_for_lateral_spread_deformation_geometry"
$ws4.Range("F5").Value = """CA landslide inventory"""
$ws4.Range("F6").Value = """none"""
$ws4.Range("D6").Value = "source for lateral spread deformation geometry: ""none"" for levels 1 and 2, and ""CPT-based"" for level 3 and directory to CPT data must be specified under ""GIS Data"" tab"
$ws4.Range("D5").Value = "source for landslide deformation geometry: ""none"" for level 1, ""CA landslide inventory"" for level 2 and 3; users can also specify path to shapefile with geometries"
$ws4.Range("B5").Value = $true
$ws4.Range("C5").Value = $true
$ws4.Range("E5").Value = "unitless"
$ws4.Range("B6").Value = $true
$ws4.Range("C6").Value = $true
$ws4.Range("E6").Value = "unitless"

# ---------------------------------------------------------------
# 4. Column widths on "fixed" sheet
# ---------------------------------------------------------------
$ws4.Columns.Item(1).ColumnWidth = 49.592447916666664
$ws4.Columns.Item(2).ColumnWidth = 16.022135416666668
$ws4.Columns.Item(3).ColumnWidth = 18.736979166666668
$ws4.Columns.Item(4).ColumnWidth = 44.877604166666664
$ws4.Columns.Item(5).ColumnWidth = 44.877604166666664
$ws4.Columns.Item(6).ColumnWidth = 44.877604166666664

# ---------------------------------------------------------------
# 5. level1!I27 value change 0.1 -> 0
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("level1")
$ws1.Range("I27").Value = 0

# ---------------------------------------------------------------
# 6. Sheet view selections (set before activating the final sheet
#    so that only the intended sheet keeps tabSelected="1")
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("level2")
$ws2.Range("C36").Select()

$ws3 = $wb.Worksheets.Item("level3")
$ws3.Range("F5").Select()

$ws4.Range("D5").Select()

$ws1.Activate()
$ws1.Range("I28").Select()
